$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for this market/product; insert
# it as a new row right before the current row 70 (which pushes the existing
# rows 70-171 down to 71-172, exactly mirroring the diff).
$ws.Rows.Item(70).Insert()

# Copy the (now-shifted) row that used to be row 70 - it is row 71 after the
# insert - into the freshly inserted blank row 70 so every column keeps the
# same constant market/product metadata.
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(70, $c).Value2 = $ws.Cells.Item(71, $c).Value2
}

# Overwrite the two columns that actually carry the new observation: the
# date (column D = 4) and the volume (column J = 10).
$ws.Cells.Item(70, 4).Value2 = 44571
$ws.Cells.Item(70, 10).Value2 = 500
